$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.489.20'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.909.82'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D5").Value = '''325.36'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = '''0.4849'
$ws.Range("E7").Value = '  +2.87%  '
$ws.Range("D8").Value = '''0.4071'
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").Value = '''0.08172'
$ws.Range("E9").Value = '  +1.89%  '
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").Value = '''23.46'
$ws.Range("E11").Value = '  +3.83%  '
$ws.Range("D12").Value = '1.944.44'
$ws.Range("E12").Value = '  +1.61%  '
$ws.Range("D13").Value = '''6.018'
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").Value = '''7.187'
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").Value = '''90.45'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = '''0.06786'
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '''1.008'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").Value = '''17.66'
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").Value = '''1.006'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '29.517.54'
$ws.Range("D22").Value = '''5.624'
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("D24").Value = '''2.190'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").Value = '2.151.37'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '''157.15'
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("D27").Value = '''6.566'
$ws.Range("E27").Value = '  +9.38%  '
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("D29").Value = '''2.118'
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").Value = '''120.46'
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("D32").Value = '''0.09534'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").Value = '''5.513'
$ws.Range("E33").Value = '  +2.96%  '
$ws.Range("D34").Value = '''3.559'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").Value = '''1.392'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").Value = '''0.02278'
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").Value = '''0.06118'
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").Value = '''1.182'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").Value = '''10.87'
$ws.Range("E39").Value = '  +8.04%  '
$ws.Range("E40").Value = '  +2.49%  '
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("D42").Value = '''0.1857'
$ws.Range("E42").Value = '  +1.30%  '
$ws.Range("D43").Value = '''1.276'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = '''2.368'
$ws.Range("E44").Value = '  -4.95%  '
$ws.Range("D45").Value = '''12.53'
$ws.Range("E45").Value = '  +3.95%  '
$ws.Range("D46").Value = '''0.07625'
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("D47").Value = '''0.5575'
$ws.Range("E47").Value = '  +1.65%  '
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("D49").Value = '''116.48'
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("D50").Value = '''72.64'
$ws.Range("E50").Value = '  +2.38%  '
$ws.Range("E51").Value = '  +2.88%  '